# "adds preliminary processor scaling data to spreadsheet"
#
# The author re-imported 'processor scaling sheet' from a fresh
# processor_scaling.txt text-query (10 rows) which replaced the older,
# longer strong_scaling.txt-sourced table (13 rows). Reproduce the
# resulting worksheet state: new/updated data in A2:C11, the trailing
# rows 12-14 removed, the B/A ratio formula re-filled down the shrunk
# range, the query-wizard's leftover cursor position, and the new
# sheet-scoped "processor_scaling" defined name that Excel creates
# for a query table's result range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("processor scaling sheet")

# Drop the three rows that no longer exist once the smaller dataset
# was re-queried (old rows 12:14 -> 50/38.56, 55/34.711, 60/38.377).
$ws.Rows("12:14").Delete()

$numProcessors = @(5, 10, 15, 20, 30, 40, 45, 50, 55, 60)
$nsPerDay      = @(4.15, 8.231, 8.955, 8.085, 7.721, 4.477, 13.879, 5.752, 9.8, 8.906)

for ($i = 0; $i -lt $numProcessors.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $numProcessors[$i]
    $ws.Cells.Item($row, 2).Value = $nsPerDay[$i]
    $ws.Cells.Item($row, 3).Formula = "=B$row/A$row"
}

# The query table wizard registers the freshly-queried range as a
# worksheet-scoped defined name.
$ws.Names.Add("processor_scaling", "='processor scaling sheet'!`$A`$2:`$B`$11")

# Leave the selection where the import wizard's dialog left it.
$ws.Activate()
$ws.Range("C15").Select()
